# Update Name of Algo
# Apply updated KNN-imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C4"   = -12.632
    "D5"   = -8.198
    "A8"   = -21.986
    "D8"   = -8.024000000000001
    "A10"  = -21.72
    "C11"  = -12.44
    "A12"  = -21.621
    "C12"  = -12.72
    "D12"  = -7.803000000000002
    "D13"  = -7.506
    "C15"  = -13.509
    "D15"  = -8.182
    "C17"  = -13.364
    "A18"  = -22.11
    "D21"  = -7.974000000000001
    "A25"  = -21.873
    "D25"  = -7.704000000000001
    "C26"  = -13.236
    "C27"  = -12.957
    "C28"  = -12.807
    "C32"  = -12.172
    "D32"  = -7.56
    "D36"  = -7.518000000000001
    "A37"  = -19.708
    "C37"  = -12.14
    "D38"  = -7.717000000000001
    "C41"  = -12.112
    "D41"  = -7.611
    "C47"  = -12.432
    "D50"  = -7.843000000000001
    "C51"  = -12.157
    "D52"  = -8.088000000000001
    "A55"  = -22.094
    "D59"  = -7.551
    "C65"  = -12.1
    "D67"  = -7.687
    "A68"  = -21.673
    "C73"  = -12.685
    "A77"  = -20.409
    "A78"  = -20.201
    "A79"  = -20.769
    "A80"  = -20.239
    "A81"  = -21.818
    "A82"  = -22.109
    "A84"  = -21.813
    "C84"  = -11.901
    "D84"  = -7.484
    "C85"  = -11.725
    "D86"  = -7.874000000000001
    "D88"  = -8.026
    "C89"  = -13.787
    "D89"  = -8.279999999999999
    "C93"  = -12.609
    "C95"  = -12.217
    "D95"  = -7.580000000000001
    "C98"  = -12.991
    "C99"  = -11.728
    "A101" = -20.326
    "C101" = -12.514
    "A102" = -20.075
    "C102" = -12.255
    "D105" = -7.777999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
